$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, pushing the existing "Holter" row down to row 6
$ws.Rows.Item(5).Insert()

# Fill the newly inserted row 5 with TTE / MR / 2 / TRUE
$ws.Range("A5").Value = "TTE"
$ws.Range("B5").Value = "MR"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = $true

# Append a new row 7 with TTE / AEKG / 2 / TRUE
$ws.Range("A7").Value = "TTE"
$ws.Range("B7").Value = "AEKG"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = $true

# Update selection to match the target state
$ws.Range("C8").Select()
